$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Month And Week")

# Fix up the fund/scheme codes that used characters which aren't safe for
# downstream systems (":", "/", "%") by swapping them for underscores
# (and "%" -> "P", e.g. "25%C" -> "25PC").
$rng = $ws.UsedRange
$rng.Replace(":", "_", -4162) | Out-Null
$rng.Replace("/", "_", -4162) | Out-Null
$rng.Replace("%", "P", -4162) | Out-Null

# Update the saved view state (scroll position / selection) on the sheet.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 95
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("F80").Select()
